$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3800153732299805
$ws.Range("E2").Value = 41.70139439737795
$ws.Range("F2").Value = 0.001419457622318641
$ws.Range("G2").Value = 0.00121615936177821
$ws.Range("H2").Value = 0.001118081889885234
$ws.Range("I2").Value = 0.001066653389605538
$ws.Range("J2").Value = 0.001048378160576709
$ws.Range("K2").Value = 0.0009785300387264544
$ws.Range("L2").Value = 0.0009502486708457281
$ws.Range("M2").Value = 0.0009502486708457281
$ws.Range("N2").Value = 0.0008956023769150044
$ws.Range("O2").Value = 0.0008956023769150044
$ws.Range("P2").Value = 0.0008631048866385359
$ws.Range("Q2").Value = 0.0008492332515253653
$ws.Range("R2").Value = 0.0008492332515253653
$ws.Range("S2").Value = 0.0008481882999216106
$ws.Range("T2").Value = 0.0008377467491839508
$ws.Range("U2").Value = 0.0008249981319727383
$ws.Range("V2").Value = 0.0008249981319727383
$ws.Range("W2").Value = 0.0008128926783114609
$ws.Range("X2").Value = 0.0008128926783114609
$ws.Range("Y2").Value = 0.0008128926783114609

$ws.Range("C3").Value = 0.408970832824707
$ws.Range("E3").Value = 42.50879382253152
$ws.Range("F3").Value = 0.001482084267898751
$ws.Range("G3").Value = 0.001205380188338145
$ws.Range("H3").Value = 0.001086320625943922
$ws.Range("I3").Value = 0.001064625485181445
$ws.Range("J3").Value = 0.001035177679680035
$ws.Range("K3").Value = 0.0009871697737845416
$ws.Range("L3").Value = 0.0009804911295697614
$ws.Range("M3").Value = 0.0009689782429288444
$ws.Range("N3").Value = 0.0009211604167707559
$ws.Range("O3").Value = 0.0009157766578668172
$ws.Range("P3").Value = 0.0009002333556871457
$ws.Range("Q3").Value = 0.0008919045015109766
$ws.Range("R3").Value = 0.000881267573710189
$ws.Range("S3").Value = 0.0008706832969711614
$ws.Range("T3").Value = 0.0008681532635271341
$ws.Range("U3").Value = 0.0008562495121876737
$ws.Range("V3").Value = 0.0008509697442718921
$ws.Range("W3").Value = 0.0008419029915435034
$ws.Range("X3").Value = 0.0008365974486657708
$ws.Range("Y3").Value = 0.0008286314585288792

$ws.Range("C4").Value = 0.3300011157989502
$ws.Range("E4").Value = 40.59347726899614
$ws.Range("F4").Value = 0.001482084267898751
$ws.Range("G4").Value = 0.001258431128012374
$ws.Range("H4").Value = 0.001094846752059726
$ws.Range("I4").Value = 0.001045014103412016
$ws.Range("J4").Value = 0.001012619121472746
$ws.Range("K4").Value = 0.001012619121472746
$ws.Range("L4").Value = 0.0009736425344691812
$ws.Range("M4").Value = 0.0009148016490252218
$ws.Range("N4").Value = 0.0008967113848651592
$ws.Range("O4").Value = 0.0008654073886145762
$ws.Range("P4").Value = 0.0008654073886145762
$ws.Range("Q4").Value = 0.0008654073886145762
$ws.Range("R4").Value = 0.0008267856954004108
$ws.Range("S4").Value = 0.0008267856954004108
$ws.Range("T4").Value = 0.0008267856954004108
$ws.Range("U4").Value = 0.0008043867667372651
$ws.Range("V4").Value = 0.0008043867667372651
$ws.Range("W4").Value = 0.0008021256370968195
$ws.Range("X4").Value = 0.0007912958531968057
$ws.Range("Y4").Value = 0.0007912958531968057

$ws.Range("C5").Value = 0.3659975528717041
$ws.Range("E5").Value = 41.78305646651825
$ws.Range("F5").Value = 0.001429509560417592
$ws.Range("G5").Value = 0.001189861407540581
$ws.Range("H5").Value = 0.00117788667309219
$ws.Range("I5").Value = 0.001106816708749912
$ws.Range("J5").Value = 0.00106068569037872
$ws.Range("K5").Value = 0.001050146826171159
$ws.Range("L5").Value = 0.0009674830778740799
$ws.Range("M5").Value = 0.0009622420407940026
$ws.Range("N5").Value = 0.0008971902869502741
$ws.Range("O5").Value = 0.0008948935950108761
$ws.Range("P5").Value = 0.0008747076175007137
$ws.Range("Q5").Value = 0.0008448053606662496
$ws.Range("R5").Value = 0.0008381547918192539
$ws.Range("S5").Value = 0.0008381547918192539
$ws.Range("T5").Value = 0.0008381547918192539
$ws.Range("U5").Value = 0.0008381547918192539
$ws.Range("V5").Value = 0.0008242546509987717
$ws.Range("W5").Value = 0.0008242546509987717
$ws.Range("X5").Value = 0.000820749128439207
$ws.Range("Y5").Value = 0.0008144845315110768

$ws.Range("C6").Value = 0.350001335144043
$ws.Range("E6").Value = 40.95489048340278
$ws.Range("F6").Value = 0.001482084267898751
$ws.Range("G6").Value = 0.001252574983202144
$ws.Range("H6").Value = 0.001111901477032265
$ws.Range("I6").Value = 0.001042187793783017
$ws.Range("J6").Value = 0.001042187793783017
$ws.Range("K6").Value = 0.001005716092924672
$ws.Range("L6").Value = 0.0009881591841851196
$ws.Range("M6").Value = 0.0009624392537924234
$ws.Range("N6").Value = 0.000934490545673592
$ws.Range("O6").Value = 0.0009192648518211696
$ws.Range("P6").Value = 0.0009114766860675195
$ws.Range("Q6").Value = 0.0008845394555320065
$ws.Range("R6").Value = 0.0008763642920680386
$ws.Range("S6").Value = 0.0008233773534034678
$ws.Range("T6").Value = 0.0008233773534034678
$ws.Range("U6").Value = 0.0008233773534034678
$ws.Range("V6").Value = 0.0008092374261901942
$ws.Range("W6").Value = 0.0008092374261901942
$ws.Range("X6").Value = 0.0007997894688520267
$ws.Range("Y6").Value = 0.0007983409450955706

$ws.Range("C7").Value = 0.4639983177185059
$ws.Range("E7").Value = 40.3026086732516
$ws.Range("F7").Value = 0.001482084267898751
$ws.Range("G7").Value = 0.001194823165296413
$ws.Range("H7").Value = 0.001108652265872566
$ws.Range("I7").Value = 0.001044114228164593
$ws.Range("J7").Value = 0.001007106897795964
$ws.Range("K7").Value = 0.0009644470256586045
$ws.Range("L7").Value = 0.0009020247088230131
$ws.Range("M7").Value = 0.0009013729203995847
$ws.Range("N7").Value = 0.0008693613866986377
$ws.Range("O7").Value = 0.0008398947473746671
$ws.Range("P7").Value = 0.0008398947473746671
$ws.Range("Q7").Value = 0.0008398947473746671
$ws.Range("R7").Value = 0.0008362965580014295
$ws.Range("S7").Value = 0.0008290638022071395
$ws.Range("T7").Value = 0.0008111969689862924
$ws.Range("U7").Value = 0.0008111969689862924
$ws.Range("V7").Value = 0.0007957125727421795
$ws.Range("W7").Value = 0.0007902290682199461
$ws.Range("X7").Value = 0.0007872350559161972
$ws.Range("Y7").Value = 0.0007856259000633837

$ws.Range("C8").Value = 0.3989934921264648
$ws.Range("E8").Value = 40.1888487761189
$ws.Range("F8").Value = 0.001414860653344698
$ws.Range("G8").Value = 0.001186955360468314
$ws.Range("H8").Value = 0.001137410728169429
$ws.Range("I8").Value = 0.001089793266306784
$ws.Range("J8").Value = 0.0009985174041095219
$ws.Range("K8").Value = 0.0009692240764308546
$ws.Range("L8").Value = 0.0009337282148615903
$ws.Range("M8").Value = 0.0009194071214166439
$ws.Range("N8").Value = 0.0009045941631004918
$ws.Range("O8").Value = 0.0009045941631004918
$ws.Range("P8").Value = 0.0008831938041804288
$ws.Range("Q8").Value = 0.0008455017208538544
$ws.Range("R8").Value = 0.0008330547659690218
$ws.Range("S8").Value = 0.0008306172793565305
$ws.Range("T8").Value = 0.0008061478350374248
$ws.Range("U8").Value = 0.0008061478350374248
$ws.Range("V8").Value = 0.0007994945910642684
$ws.Range("W8").Value = 0.000791767331614247
$ws.Range("X8").Value = 0.0007834083582089453
$ws.Range("Y8").Value = 0.0007834083582089453

$ws.Range("C9").Value = 0.4199998378753662
$ws.Range("E9").Value = 41.23531659502623
$ws.Range("F9").Value = 0.001482084267898751
$ws.Range("G9").Value = 0.001227186079113327
$ws.Range("H9").Value = 0.001133533785775148
$ws.Range("I9").Value = 0.001114257433058895
$ws.Range("J9").Value = 0.00106096290803818
$ws.Range("K9").Value = 0.0009760901623953594
$ws.Range("L9").Value = 0.0009760901623953594
$ws.Range("M9").Value = 0.0009155600396585376
$ws.Range("N9").Value = 0.0009155600396585376
$ws.Range("O9").Value = 0.0009155600396585376
$ws.Range("P9").Value = 0.0008943389402854208
$ws.Range("Q9").Value = 0.0008630139930225026
$ws.Range("R9").Value = 0.0008627184516832593
$ws.Range("S9").Value = 0.0008454026360785587
$ws.Range("T9").Value = 0.000840703418672778
$ws.Range("U9").Value = 0.0008399777871298346
$ws.Range("V9").Value = 0.0008256173065408226
$ws.Range("W9").Value = 0.0008118286973277226
$ws.Range("X9").Value = 0.000810588072812532
$ws.Range("Y9").Value = 0.0008038073410336496

$ws.Range("C10").Value = 0.349998950958252
$ws.Range("E10").Value = 40.32145678267443
$ws.Range("F10").Value = 0.001458791654765865
$ws.Range("G10").Value = 0.001216276843656992
$ws.Range("H10").Value = 0.00113466043471599
$ws.Range("I10").Value = 0.001061990265907244
$ws.Range("J10").Value = 0.001029809739242654
$ws.Range("K10").Value = 0.0009916263387956283
$ws.Range("L10").Value = 0.0009445646693569918
$ws.Range("M10").Value = 0.0009423668477601712
$ws.Range("N10").Value = 0.0009080063401603752
$ws.Range("O10").Value = 0.000903424850859583
$ws.Range("P10").Value = 0.0008651260973563597
$ws.Range("Q10").Value = 0.0008651260973563597
$ws.Range("R10").Value = 0.0008356869955482675
$ws.Range("S10").Value = 0.000828295973822084
$ws.Range("T10").Value = 0.0008054972648465141
$ws.Range("U10").Value = 0.0008054972648465141
$ws.Range("V10").Value = 0.0008054972648465141
$ws.Range("W10").Value = 0.0007941633791962838
$ws.Range("X10").Value = 0.0007907142778933145
$ws.Range("Y10").Value = 0.0007859933096037899

$ws.Range("C11").Value = 0.3409979343414307
$ws.Range("E11").Value = 42.12382639776479
$ws.Range("F11").Value = 0.001482084267898751
$ws.Range("G11").Value = 0.001225818594533986
$ws.Range("H11").Value = 0.001107199958339983
$ws.Range("I11").Value = 0.0009716175156113394
$ws.Range("J11").Value = 0.0009374628563379308
$ws.Range("K11").Value = 0.0009374628563379308
$ws.Range("L11").Value = 0.0009374628563379308
$ws.Range("M11").Value = 0.0009298460148008042
$ws.Range("N11").Value = 0.0009285560659384704
$ws.Range("O11").Value = 0.0008847614379854943
$ws.Range("P11").Value = 0.0008847614379854943
$ws.Range("Q11").Value = 0.0008647706536492328
$ws.Range("R11").Value = 0.0008647706536492328
$ws.Range("S11").Value = 0.0008636799663276758
$ws.Range("T11").Value = 0.0008484660196850473
$ws.Range("U11").Value = 0.0008382543265603341
$ws.Range("V11").Value = 0.00083158685391789
$ws.Range("W11").Value = 0.0008280542560920829
$ws.Range("X11").Value = 0.0008211272202293329
$ws.Range("Y11").Value = 0.0008211272202293329
